$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching style of existing headers (e.g. G1)
# by copying G1's formatting (keeps the same shared cell style index).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the numeric values for the new Save column
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
